$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Moon Dates")

# Row 14 held a duplicate full-moon reading (12/31/2009, Full, Capricorn --
# a repeat of the 1/11/2009 row). Clear its contents: the row stays in place
# but A14:C14 become blank (A14 keeps its date-cell formatting).
$ws.Range("A14:C14").ClearContents()

# The sheet also had a spare trailing blank row (26) past the real data.
# Remove it outright so the used range shrinks from A1:C26 to A1:C25.
$ws.Rows.Item(26).Delete()

# Re-apply the AutoFilter over the smaller range so its stored reference
# shrinks too (A1:C26 -> A1:C25) instead of pointing past the data.
$ws.AutoFilterMode = $false
$ws.Range("A1:C25").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name (which mirrors the AutoFilter
# range) in sync with the new extent.
$wb.Names.Item($ws.Name + "!_FilterDatabase").RefersTo = "='" + $ws.Name + "'!`$A`$1:`$C`$25"

# Leave the selection where the user ended up after the edit.
$ws.Range("D7").Select()

$wb.Save()
